# Applies the changes described by the diff:
#  - Rename sheet "sheet_cut_off" (1st tab) -> "CUT OFFS"
#  - Rename sheet "Sheet1" (2nd tab) -> "sheet_cut_off"
#  - On the renamed "sheet_cut_off" (formerly Sheet1) sheet, change A2:A6 text labels
#    to the new "_score" names
#  - Make the (new) "sheet_cut_off" sheet the active/selected tab,
#    and select cell G17 on it
#  - Leave the "CUT OFFS" sheet not selected/active

$wb = $excel.ActiveWorkbook

$wsCutOffs = $wb.Worksheets.Item(1)   # was "sheet_cut_off", becomes "CUT OFFS"
$wsSheet1  = $wb.Worksheets.Item(2)   # was "Sheet1", becomes "sheet_cut_off"

# Update the cell labels on the second worksheet before renaming (values are
# independent of the sheet name, but do it in a sensible order anyway).
$wsSheet1.Range("A2").Value = "immunity_score"
$wsSheet1.Range("A3").Value = "survaillance_score"
$wsSheet1.Range("A4").Value = "determinant_score"
$wsSheet1.Range("A5").Value = "outbreak_score"
$wsSheet1.Range("A6").Value = "total_score"

# Rename the sheets.
$wsCutOffs.Name = "CUT OFFS"
$wsSheet1.Name = "sheet_cut_off"

# Activate the "sheet_cut_off" sheet (formerly Sheet1) and select G17 on it.
$wsSheet1.Activate()
$wsSheet1.Range("G17").Select()
